$d = $word.ActiveDocument

function Replace-Text($find, $replace, $wholeWord) {
    $result = $d.Content.Find.Execute($find, $false, $wholeWord, $false, $false, $false,
                                       $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Output ("WARNING: not found -> " + $find)
    }
}

# Receipt / work-order number ("ÁTVÉELI ELISMERVÉNY ... 14" -> "... 37")
Replace-Text "14" "37" $true

# Customer address
Replace-Text "Cím: 9700 Szombathely Nemlétezik utca. 3" "Cím: 9700 Szombathely Nincs Ilyen út 69" $false

# Comment field
Replace-Text "Megjegyzés:uuheuhehu" "Megjegyzés:Valamien megjegyzés" $false

# Product details
Replace-Text "Megnevezés: Fünyiro" "Megnevezés: Fûnyíró" $false
Replace-Text "Típus: Ferrari" "Típus: Husqwarna" $false
Replace-Text "Modell: Igen" "Modell: CW23" $false

# Fault description / accessories
Replace-Text "Hibajelenség: Van" "Hibajelenség: Nem indul el" $false
Replace-Text "Tartozékok: Nincs" "Tartozékok: Kerék, Kesztyû" $false

# Service diagnosis
Replace-Text "Szerviz diagnózis: Rósz" "Szerviz diagnózis: El van törve" $false

# Dates (three occurrences in the document, all the same old/new value)
Replace-Text "2024.07.28" "2024.08.10" $false

Write-Output "Done"
